# Apply the edit described by the diff: add a "Created" marker to column F
# (LinkedIn_Poster) of row 2, and update the active selection, as happens
# when the file is opened, edited, and re-saved in real Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value - this adds a new shared string "Created" and writes it
# into cell F2.
$ws.Range("F2").Value = "Created"

# Reflect the updated selection seen in the saved file.
$ws.Range("G7").Select()

$wb.Save()
